$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - (Intercept)
$ws.Range("B2").Value = 2593.792633
$ws.Range("D2").Value = 44.350781

# Row 3 - household_group_collapsed
$ws.Range("B3").Value = 1121.284486
$ws.Range("D3").Value = 9.586318
$ws.Range("E3").Value = 0.000101

# Row 4 - Residuals
$ws.Range("B4").Value = 13041.839241
$ws.Range("C4").Value = 223

# Row 5 - SM-Control
$ws.Range("G5").Value = -3.07613
$ws.Range("H5").Value = -6.41605
$ws.Range("I5").Value = 0.26379
$ws.Range("J5").Value = 0.078071

# Row 6 - SM + Traps-Control
$ws.Range("G6").Value = 1.863805
$ws.Range("H6").Value = -1.732582
$ws.Range("I6").Value = 5.460192
$ws.Range("J6").Value = 0.440974

# Row 7 - SM + Traps-SM
$ws.Range("G7").Value = 4.939935
$ws.Range("H7").Value = 2.221029
$ws.Range("I7").Value = 7.658841
$ws.Range("J7").Value = 0.00008
